$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Day 1 of the schedule had a duplicate "Tutorial" entry (rows 9 and 15).
# Remove the first occurrence (row 9); this shifts all subsequent rows up by one.
$ws.Rows.Item(9).Delete()

# The remaining "Tutorial" entry (now on row 14 after the shift) is renamed
# to "Indexing" to reflect the newly added indexing topic for day 1.
$ws.Range("B14").Value = "Indexing"

# Update the active cell selection to match the edited cell.
$ws.Range("B14").Select()
